$d = $word.ActiveDocument

$pairs = @(
    @("2025-10-27 Monday", "2025-10-28 Tuesday"),
    @("55×59=", "86×54="),
    @("52×47=", "12×41="),
    @("33×53=", "16×53="),
    @("75×90=", "25×23="),
    @("37×50=", "73×47="),
    @("94×49=", "86×71="),
    @("22×98=", "15×60="),
    @("28×68=", "27×61="),
    @("40×38=", "12×85="),
    @("96×21=", "82×51="),
    @("19×22=", "53×15="),
    @("94×34=", "71×55="),
    @("52×48=", "69×58="),
    @("82×95=", "54×16="),
    @("68×28=", "26×85="),
    @("13×25=", "29×83="),
    @("72×41=", "75×26="),
    @("15×31=", "96×59="),
    @("90×45=", "56×76="),
    @("94×86=", "84×46="),
    @("32×98=", "55×57="),
    @("55×26=", "69×69="),
    @("53×73=", "27×91="),
    @("36×44=", "32×17="),
    @("25×60=", "40×20=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
